$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Principal")
$ws2 = $wb.Worksheets.Item("Doctores")
$ws3 = $wb.Worksheets.Item("Pacientes")

# ---- Sheet "Principal": fill in row 3 with the alternate ("error") case ----
$ws1.Range("D3").Value = "asanchez"
$ws1.Range("G3").Value = "Cita Error"
$ws1.Range("E3").Value = "juribe"

# ---- Sheet "Pacientes": add the new patient record used by the new case ----
$ws3.Range("C4").Value = "Uribe"
$ws3.Range("A4").Value = "juribe"
$ws3.Range("B4").Value = "Juan D."
$ws3.Range("E4").Value = 1010198
$ws3.Range("F4").Value = 10001009

$ws1.Range("F3").Value = "20/06/2018"
$ws1.Range("H3").Value = "Error:"

# ---- Sheet "Doctores": update the Documento value ----
$ws2.Range("F4").Value = 10299990

# ---- restore selections to match the edited cells ----
$ws2.Range("F4").Select()
$ws3.Range("A9").Select()
$ws1.Select()
$ws1.Range("F3").Select()
